# Actualización automática del mapa (2025-09-03 22:05:55)
# Agrega dos nuevas filas (91 y 92) a la hoja PEBCOM, reproduciendo el mismo
# formato que el resto de los datos: todas las columnas se guardan como texto
# excepto I (Attachments), M (Coordenada_X) y N (Coordenada_Y), que son numéricas.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($Worksheet, [string]$Address, [string]$Text)
    $cell = $Worksheet.Range($Address)
    $cell.NumberFormat = "@"
    $cell.Value = $Text
    $cell.Style = "Normal"
}

function Set-NumberCell {
    param($Worksheet, [string]$Address, [double]$Number)
    $Worksheet.Range($Address).Value = $Number
}

$newRows = @(
    @{
        Row = 91
        A = "7136"
        B = "9/3/2025"
        C = "FERRARI 455"
        D = "15"
        E = "809427020"
        F = "PEBCOM"
        G = "Pendiente"
        H = "Picada"
        I = 1
        J = "Cambio"
        K = "Sin equipos"
        L = "Terminal"
        M = -58.441587
        N = -34.60547
        O = "Paternal"
        P = "Capital Norte"
    },
    @{
        Row = 92
        A = "7150"
        B = "9/3/2025"
        C = "Bartolomé Mitre 3070"
        D = "3"
        E = "809427021"
        F = "PEBCOM"
        G = "Pendiente"
        H = "Cambiar"
        I = 1
        J = "Cambio"
        K = "Sin equipos"
        L = "Pasante"
        M = -58.410025
        N = -34.609184
        O = "Almagro"
        P = "Capital Sur"
    }
)

foreach ($r in $newRows) {
    $row = $r.Row
    Set-TextCell $ws "A$row" $r.A
    Set-TextCell $ws "B$row" $r.B
    Set-TextCell $ws "C$row" $r.C
    Set-TextCell $ws "D$row" $r.D
    Set-TextCell $ws "E$row" $r.E
    Set-TextCell $ws "F$row" $r.F
    Set-TextCell $ws "G$row" $r.G
    Set-TextCell $ws "H$row" $r.H
    Set-NumberCell $ws "I$row" $r.I
    Set-TextCell $ws "J$row" $r.J
    Set-TextCell $ws "K$row" $r.K
    Set-TextCell $ws "L$row" $r.L
    Set-NumberCell $ws "M$row" $r.M
    Set-NumberCell $ws "N$row" $r.N
    Set-TextCell $ws "O$row" $r.O
    Set-TextCell $ws "P$row" $r.P
}
